$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MP_PFOA_PFOS")

# Fill column A (Replicates) for rows 2-19 with a repeating 1,2,3 pattern
for ($row = 2; $row -le 19; $row++) {
    $value = (($row - 2) % 3) + 1
    $ws.Cells.Item($row, 1).Value = $value
}

# Activate the sheet and update the current selection to match the saved view
$ws.Activate()
$ws.Range("A17:A19").Select()
